# Scheduled runner refresh: update computed price/profit columns (H-N) across sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 216
$ws.Range("I5").Value = 216
$ws.Range("K5").Value = 216
$ws.Range("M5").Value = -101
$ws.Range("H12").Value = 517.4
$ws.Range("J12").Value = 517.4
$ws.Range("L12").Value = 517.4
$ws.Range("N12").Value = -857.4
$ws.Range("H17").Value = 1342.2858
$ws.Range("J17").Value = 1432.6666
$ws.Range("L17").Value = 4297.9998
$ws.Range("N17").Value = -4633.9998
$ws.Range("H19").Value = 764
$ws.Range("I19").Value = 574.5
$ws.Range("J19").Value = 839.8
$ws.Range("K19").Value = 574.5
$ws.Range("L19").Value = 839.8
$ws.Range("M19").Value = -399.5
$ws.Range("N19").Value = -1189.8
$ws.Range("H33").Value = 315.63635
$ws.Range("I33").Value = 92.47619
$ws.Range("K33").Value = 92.47619
$ws.Range("M33").Value = 136.52381
$ws.Range("H58").Value = 150
$ws.Range("J58").Value = 150
$ws.Range("L58").Value = 450
$ws.Range("N58").Value = -750
$ws.Range("H106").Value = 2473.75
$ws.Range("I106").Value = 2298.3333
$ws.Range("K106").Value = 2298.3333
$ws.Range("M106").Value = -1667.3333
$ws.Range("H129").Value = 1799
$ws.Range("I129").Value = 1799
$ws.Range("K129").Value = 5397
$ws.Range("M129").Value = -397
$ws.Range("H135").Value = 799.0769
$ws.Range("I135").Value = 761.5
$ws.Range("K135").Value = 6853.5
$ws.Range("M135").Value = -4318.5
$ws.Range("H137").Value = 2099.8
$ws.Range("J137").Value = 2249.5
$ws.Range("L137").Value = 6748.5
$ws.Range("N137").Value = -11848.5
$ws.Range("H138").Value = 3626.1667
$ws.Range("I138").Value = 2350
$ws.Range("J138").Value = 3881.4
$ws.Range("K138").Value = 7050
$ws.Range("L138").Value = 11644.2
$ws.Range("M138").Value = -1910
$ws.Range("N138").Value = -21924.2

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 175
$ws.Range("I5").Value = 175
$ws.Range("K5").Value = 175
$ws.Range("M5").Value = -63
$ws.Range("H32").Value = 2099.625
$ws.Range("I32").Value = 2113.9285
$ws.Range("K32").Value = 2113.9285
$ws.Range("M32").Value = -1826.9285
$ws.Range("H95").Value = 671999.7
$ws.Range("J95").Value = 671999.7
$ws.Range("L95").Value = 671999.7
$ws.Range("N95").Value = -677491.7
$ws.Range("H101").Value = 75000
$ws.Range("J101").Value = 75000
$ws.Range("L101").Value = 75000
$ws.Range("N101").Value = -81490
$ws.Range("H111").Value = 0
$ws.Range("J111").Value = 0
$ws.Range("L111").Value = 0
$ws.Range("H125").Value = 100000
$ws.Range("J125").Value = 100000
$ws.Range("L125").Value = 100000
$ws.Range("N125").Value = -109840

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 175
$ws.Range("I4").Value = 175
$ws.Range("K4").Value = 175
$ws.Range("M4").Value = -60
$ws.Range("H20").Value = 3889.3333
$ws.Range("I20").Value = 3889.3333
$ws.Range("K20").Value = 3889.3333
$ws.Range("M20").Value = -3642.3333
$ws.Range("H82").Value = 35128.5
$ws.Range("I82").Value = 35128.5
$ws.Range("K82").Value = 35128.5
$ws.Range("M82").Value = -34745.5
$ws.Range("H85").Value = 35128.5
$ws.Range("I85").Value = 35128.5
$ws.Range("K85").Value = 35128.5
$ws.Range("M85").Value = -33802.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 186.66667
$ws.Range("I22").Value = 154.28572
$ws.Range("K22").Value = 154.28572
$ws.Range("M22").Value = 195.71428
$ws.Range("H62").Value = 12500
$ws.Range("I62").Value = 12500
$ws.Range("K62").Value = 12500
$ws.Range("M62").Value = -11876
$ws.Range("H65").Value = 12500
$ws.Range("I65").Value = 12500
$ws.Range("K65").Value = 62500
$ws.Range("M65").Value = -59380
$ws.Range("H132").Value = 3056.111
$ws.Range("I132").Value = 2851
$ws.Range("J132").Value = 3466.3333
$ws.Range("K132").Value = 8553
$ws.Range("L132").Value = 10398.9999
$ws.Range("M132").Value = -6023
$ws.Range("N132").Value = -15458.9999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 36.375
$ws.Range("I7").Value = 52
$ws.Range("J7").Value = 10.333333
$ws.Range("K7").Value = 156
$ws.Range("L7").Value = 30.999999
$ws.Range("M7").Value = -44
$ws.Range("N7").Value = -254.999999
$ws.Range("H23").Value = 636.25
$ws.Range("I23").Value = 358.33334
$ws.Range("J23").Value = 803
$ws.Range("K23").Value = 1075.00002
$ws.Range("L23").Value = 2409
$ws.Range("M23").Value = -840.0000199999999
$ws.Range("N23").Value = -2879
$ws.Range("H38").Value = 70.22221999999999
$ws.Range("J38").Value = 96.666664
$ws.Range("L38").Value = 289.999992
$ws.Range("N38").Value = -983.999992
$ws.Range("H69").Value = 22269
$ws.Range("I69").Value = 24358.334
$ws.Range("J69").Value = 16001
$ws.Range("K69").Value = 73075.00199999999
$ws.Range("L69").Value = 48003
$ws.Range("M69").Value = -72264.00199999999
$ws.Range("N69").Value = -49625
$ws.Range("H72").Value = 22269
$ws.Range("I72").Value = 24358.334
$ws.Range("J72").Value = 16001
$ws.Range("K72").Value = 219225.006
$ws.Range("L72").Value = 144009
$ws.Range("M72").Value = -215169.006
$ws.Range("N72").Value = -152121
$ws.Range("H76").Value = 1000
$ws.Range("I76").Value = 1000
$ws.Range("K76").Value = 3000
$ws.Range("M76").Value = -2617
$ws.Range("H79").Value = 1000
$ws.Range("I79").Value = 1000
$ws.Range("K79").Value = 3000
$ws.Range("M79").Value = -1674
$ws.Range("H80").Value = 6000
$ws.Range("I80").Value = 4000
$ws.Range("K80").Value = 12000
$ws.Range("M80").Value = -11064
$ws.Range("H83").Value = 6000
$ws.Range("I83").Value = 4000
$ws.Range("K83").Value = 36000
$ws.Range("M83").Value = -31320
$ws.Range("H97").Value = 764.8333
$ws.Range("I97").Value = 717
$ws.Range("K97").Value = 2151
$ws.Range("M97").Value = -1655

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H20").Value = 14376.375
$ws.Range("I20").Value = 13168.5
$ws.Range("J20").Value = 18000
$ws.Range("K20").Value = 13168.5
$ws.Range("L20").Value = 18000
$ws.Range("M20").Value = -12942.5
$ws.Range("N20").Value = -18452
$ws.Range("H22").Value = 876.6786
$ws.Range("I22").Value = 967.7059
$ws.Range("J22").Value = 736
$ws.Range("K22").Value = 967.7059
$ws.Range("L22").Value = 736
$ws.Range("M22").Value = -672.7059
$ws.Range("N22").Value = -1326
$ws.Range("H27").Value = 876.6786
$ws.Range("I27").Value = 967.7059
$ws.Range("J27").Value = 736
$ws.Range("K27").Value = 967.7059
$ws.Range("L27").Value = 736
$ws.Range("M27").Value = -860.7059
$ws.Range("N27").Value = -950
$ws.Range("H132").Value = 6873.5
$ws.Range("I132").Value = 6873.5
$ws.Range("K132").Value = 20620.5
$ws.Range("M132").Value = -18090.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H31").Value = 30017
$ws.Range("I31").Value = 30017
$ws.Range("K31").Value = 30017
$ws.Range("M31").Value = -29669
$ws.Range("H132").Value = 2683.375
$ws.Range("I132").Value = 1994.5
$ws.Range("K132").Value = 5983.5
$ws.Range("M132").Value = -3453.5

# Cells whose derived value no longer applies this run are cleared entirely.
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("N111").ClearContents()
